$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.623.43"
$ws.Range("E2").Value = "  +4.24%  "
$ws.Range("D3").Value = "3.255.46"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'578.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").Value = "'181.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.20%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.90%  "
$ws.Range("D9").Value = "3.250.08"
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("D11").Value = "'6.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("D12").Value = "'0.415"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.75%  "
$ws.Range("D13").Value = "3.810.36"
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "'28.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "67.552.56"
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("D18").Value = "3.251.32"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").Value = "'5.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("D20").Value = "'13.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.19%  "
$ws.Range("D21").Value = "'375.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.38%  "
$ws.Range("D22").Value = "'7.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.07%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'71.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.511"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("D26").Value = "'0.0000120"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").Value = "'9.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'0.181"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'5.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.28%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("D32").Value = "'22.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'1.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.82%  "
$ws.Range("D35").Value = "'6.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("D36").Value = "'162.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.67%  "
$ws.Range("E37").Value = "  +3.13%  "
$ws.Range("D38").Value = "'0.855"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("D39").Value = "'1.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.64%  "
$ws.Range("D40").Value = "'6.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.94%  "
$ws.Range("D41").Value = "'26.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.71%  "
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("D43").Value = "'362.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.25%  "
$ws.Range("E44").Value = "  +6.38%  "
$ws.Range("D45").Value = "2.745.00"
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("D46").Value = "'25.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.89%  "
$ws.Range("D47").Value = "'40.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("D48").Value = "'0.0675"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.44%  "
$ws.Range("E51").Value = "  -0.46%  "
